# Thêm chức năng Import Student
# Consolidate the "NIEN KHOA" (academic-year) column (J) on the
# IMPORT_TEMPLATE sheet: every data row (2..93) currently holds a unique
# "20xx - 2019" label; replace them all with the single label "K15".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, "J").End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 93 }

for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 10).Value = "K15"
}

# Restore the scrolled/selected view the author left the sheet in.
$ws.Range("J6").Select()
$excel.ActiveWindow.ScrollColumn = 3
